$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://ncpi-fhir.github.io/ncpi-fhir-ig-v02/ValueSet/condition-inheritance-vs"
$wsMeta.Range("B3").Value = "0.2.0"
$wsMeta.Range("B8").Value = "2022-05-26T18:07:50+00:00"

$wsCodes = $wb.Worksheets.Item("Include from Condition Inheri")
$wsCodes.Range("B4").Value = "https://ncpi-fhir.github.io/ncpi-fhir-ig-v02/CodeSystem/ConditionInheritanceMode"
